$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.13692466666667
$ws.Range("H2").Value = 51.410774
$ws.Range("I2").Value = 0.538539545062134
$ws.Range("J2").Value = 0.538539545062134
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.376025
$ws.Range("N2").Value = 4.128075
$ws.Range("O2").Value = 0.1961144671983135
$ws.Range("P2").Value = 0.1961144671983135
$ws.Range("Q2").Value = 23.58083676445001
$ws.Range("R2").Value = 212.22753088005
$ws.Range("S2").Value = 0.1056153959450825
$ws.Range("T2").Value = 0.1056153959450825

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.13692466666667
$ws.Range("H3").Value = 51.410774
$ws.Range("I3").Value = 0.538539545062134
$ws.Range("J3").Value = 0.538539545062134
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.250698333333333
$ws.Range("N3").Value = 6.752095
$ws.Range("O3").Value = 0.3207750618381198
$ws.Range("P3").Value = 0.3207750618381198
$ws.Range("Q3").Value = 38.57004778572556
$ws.Range("R3").Value = 347.13043007153
$ws.Range("S3").Value = 0.1727500558695789
$ws.Range("T3").Value = 0.1727500558695789

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.13692466666667
$ws.Range("H4").Value = 51.410774
$ws.Range("I4").Value = 0.538539545062134
$ws.Range("J4").Value = 0.538539545062134
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.389714666666666
$ws.Range("N4").Value = 10.169144
$ws.Range("O4").Value = 0.4831104709635667
$ws.Range("P4").Value = 0.4831104709635668
$ws.Range("Q4").Value = 58.08928488416178
$ws.Range("R4").Value = 522.803563957456
$ws.Range("S4").Value = 0.2601740932474725
$ws.Range("T4").Value = 0.2601740932474725

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.55317066666667
$ws.Range("H5").Value = 34.659512
$ws.Range("I5").Value = 0.3630662674822902
$ws.Range("J5").Value = 0.3630662674822902
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.376025
$ws.Range("N5").Value = 4.128075
$ws.Range("O5").Value = 0.1961144671983135
$ws.Range("P5").Value = 0.1961144671983135
$ws.Range("Q5").Value = 15.8974516666
$ws.Range("R5").Value = 143.0770649994
$ws.Range("S5").Value = 0.07120254760496969
$ws.Range("T5").Value = 0.0712025476049697

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.55317066666667
$ws.Range("H6").Value = 34.659512
$ws.Range("I6").Value = 0.3630662674822902
$ws.Range("J6").Value = 0.3630662674822902
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.250698333333333
$ws.Range("N6").Value = 6.752095
$ws.Range("O6").Value = 0.3207750618381198
$ws.Range("P6").Value = 0.3207750618381198
$ws.Range("Q6").Value = 26.00270196418222
$ws.Range("R6").Value = 234.02431767764
$ws.Range("S6").Value = 0.1164626044029669
$ws.Range("T6").Value = 0.116462604402967

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.55317066666667
$ws.Range("H7").Value = 34.659512
$ws.Range("I7").Value = 0.3630662674822902
$ws.Range("J7").Value = 0.3630662674822902
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.389714666666666
$ws.Range("N7").Value = 10.169144
$ws.Range("O7").Value = 0.4831104709635667
$ws.Range("P7").Value = 0.4831104709635668
$ws.Range("Q7").Value = 39.16195205530311
$ws.Range("R7").Value = 352.4575684977279
$ws.Range("S7").Value = 0.1754011154743535
$ws.Range("T7").Value = 0.1754011154743535

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.131012
$ws.Range("H8").Value = 9.393036
$ws.Range("I8").Value = 0.09839418745557586
$ws.Range("J8").Value = 0.09839418745557586
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.376025
$ws.Range("N8").Value = 4.128075
$ws.Range("O8").Value = 0.1961144671983135
$ws.Range("P8").Value = 0.1961144671983135
$ws.Range("Q8").Value = 4.3083507873
$ws.Range("R8").Value = 38.7751570857
$ws.Range("S8").Value = 0.01929652364826124
$ws.Range("T8").Value = 0.01929652364826124

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.131012
$ws.Range("H9").Value = 9.393036
$ws.Range("I9").Value = 0.09839418745557586
$ws.Range("J9").Value = 0.09839418745557586
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.250698333333333
$ws.Range("N9").Value = 6.752095
$ws.Range("O9").Value = 0.3207750618381198
$ws.Range("P9").Value = 0.3207750618381198
$ws.Range("Q9").Value = 7.046963490046667
$ws.Range("R9").Value = 63.42267141042
$ws.Range("S9").Value = 0.0315624015655739
$ws.Range("T9").Value = 0.0315624015655739

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.131012
$ws.Range("H10").Value = 9.393036
$ws.Range("I10").Value = 0.09839418745557586
$ws.Range("J10").Value = 0.09839418745557586
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.389714666666666
$ws.Range("N10").Value = 10.169144
$ws.Range("O10").Value = 0.4831104709635667
$ws.Range("P10").Value = 0.4831104709635668
$ws.Range("Q10").Value = 10.61323729790933
$ws.Range("R10").Value = 95.519135681184
$ws.Range("S10").Value = 0.04753526224174073
$ws.Range("T10").Value = 0.04753526224174073
